# Remove the "hk02" row (row 3) from the scattering-planes table.
# Deleting the entire row shifts every row below it up by one, which
# matches the target sheet (dimension shrinks from A1:J38 to A1:J37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Delete()
